# Disaggregation of commodity Copper
#
# 1) Rename the "Copper ores and concentrates" category to "Copper" on every
#    yearly sheet (column C, row 7 of each sheet holds that category label as
#    a shared string — rewriting it on every sheet lets the engine drop the
#    now-unused old string and converge all the cells on the new text).
#
# 2) Re-arrange the disaggregated recycling figures: on every sheet the
#    numbers that used to sit in columns D/E/F of rows 5-8 are rotated one
#    column to the right (D->E, E->F, F->D) to reflect the new commodity
#    breakdown introduced by splitting "Copper ores and concentrates" into
#    "Copper".

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Rename the category label in column C, row 7 on every sheet.
    $ws.Cells.Item(7, 3).Value2 = "Copper"

    # Rotate columns D, E, F (4, 5, 6) to the right for rows 5 through 8.
    for ($r = 5; $r -le 8; $r++) {
        $dVal = $ws.Cells.Item($r, 4).Value2
        $eVal = $ws.Cells.Item($r, 5).Value2
        $fVal = $ws.Cells.Item($r, 6).Value2

        $ws.Cells.Item($r, 4).Value2 = $fVal
        $ws.Cells.Item($r, 5).Value2 = $dVal
        $ws.Cells.Item($r, 6).Value2 = $eVal
    }
}
